$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new LoginServer row (row 2) with server data.
# Write order chosen to match the shared-string insertion order:
# 127.0.0.1, LoginServer_1, 000106001
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "LoginServer_1"
$ws.Range("C2").Value = "LoginServer_1"
$ws.Range("B2").Value = "000106001"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 6001

# A2/B2 already carried the text number format; extend it to C2/F2 as well.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

# Move the selection to G5 (was A2:H7 with active cell H7).
$ws.Range("G5").Select()

# The TRUE/FALSE list validation now starts at row 3 instead of row 2.
$ws.Range("F2").Validation.Delete()
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
